$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 14
$wsExhibit.Range("F4").Value = 501

# Sheet "全部类型" (all types) - mirrors the same data
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 14
$wsAll.Range("F4").Value = 501
